$d = $word.ActiveDocument

# The document currently ends with a paragraph styled "Bibliography" that
# contains just a period, followed by the "_GoBack" bookmark.
# We will: (1) split that paragraph so a brand-new, unstyled paragraph is
# inserted just before it, (2) fill that new paragraph with the PostgreSQL
# sentence, (3) move the _GoBack bookmark onto the end of that new
# paragraph, and (4) turn the old "." paragraph into an empty, justified
# paragraph.

$pBib = $d.Paragraphs.Item($d.Paragraphs.Count)
$splitPoint = $d.Range($pBib.Range.Start, $pBib.Range.Start)
$splitPoint.InsertParagraphBefore()

$pNew = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$pNew.Style = $d.Styles.Item("Normal")
$pNew.Format.Hyphenation = -1
$pNew.Range.Text = "After the difficulties I observed with SQLite, I migrate the experiment to PostgreSQL, which does implement merge-join and hash-join. "

# Remove the old bookmark and recreate it at the end of the new paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$pNewAfterText = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$bmRange = $d.Range($pNewAfterText.Range.End - 1, $pNewAfterText.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Clear out the trailing "." paragraph and turn it into an empty, justified
# paragraph.
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$textOnly = $d.Range($pLast.Range.Start, $pLast.Range.End - 1)
$textOnly.Text = ""
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$pLast.Style = $d.Styles.Item("Normal")
$pLast.Format.Hyphenation = -1
$pLast.Alignment = 3
